$wb = $excel.ActiveWorkbook

# --- 1. Update the existing "FiP" sheet's view: activate it first so we can
#        move its selection, then it will be superseded as "active" once the
#        new sheet is activated below (matches the diff: tabSelected removed,
#        selection moved from N7 to F23). ---
$fip = $wb.Worksheets.Item("FiP")
$fip.Activate()
$fip.Range("F23").Select()

# --- 2. Add the new "SiP" (Sixth Preference) sheet after "FiP" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SiP"

# --- Row 1: headers ---
$ws.Range("A1").Value = "No. of responses"
$ws.Range("B1").Value = "DIP"
$ws.Range("C1").Value = "AIS"
$ws.Range("D1").Value = "SEO"
$ws.Range("E1").Value = "SA"
$ws.Range("F1").Value = "UE"
$ws.Range("G1").Value = "ACN"
$ws.Range("H1").Value = "Sum"
$ws.Range("J1").Value = "DIP"
$ws.Range("K1").Value = "AIS"
$ws.Range("L1").Value = "SEO"
$ws.Range("M1").Value = "SA"
$ws.Range("N1").Value = "UE"
$ws.Range("O1").Value = "ACN"
$ws.Range("P1").Value = "Total"

# --- Row 2 ---
$ws.Range("A2").Value = 19
$ws.Range("B2").Formula = "=A2*J2%"
$ws.Range("C2").Formula = "=A2*K2%"
$ws.Range("D2").Formula = "=A2*L2%"
$ws.Range("E2").Formula = "=A2*M2%"
$ws.Range("F2").Formula = "=A2*N2%"
$ws.Range("G2").Formula = "=A2*O2%"
$ws.Range("H2").Formula = "=SUM(B2:G2)"
$ws.Range("J2").Value = 26.3
$ws.Range("K2").Value = 10.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 21.1
$ws.Range("N2").Value = 5.3
$ws.Range("O2").Value = 36.8
$ws.Range("P2").Formula = "=SUM(J2:O2)"

# --- Row 3 ---
$ws.Range("A3").Value = 38
$ws.Range("B3").Formula = "=A3*J3%"
$ws.Range("C3").Formula = "=A3*K3%"
$ws.Range("D3").Formula = "=A3*L3%"
$ws.Range("E3").Formula = "=A3*M3%"
$ws.Range("F3").Formula = "=A3*N3%"
$ws.Range("G3").Formula = "=A3*O3%"
$ws.Range("H3").Formula = "=SUM(B3:G3)"
$ws.Range("J3").Value = 31.6
$ws.Range("K3").Value = 5.3
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 18.4
$ws.Range("N3").Value = 5.3
$ws.Range("O3").Value = 39.5
$ws.Range("P3").Formula = "=SUM(J3:O3)"

# --- Row 4 ---
$ws.Range("A4").Value = 54
$ws.Range("B4").Formula = "=A4*J4%"
$ws.Range("C4").Formula = "=A4*K4%"
$ws.Range("D4").Formula = "=A4*L4%"
$ws.Range("E4").Formula = "=A4*M4%"
$ws.Range("F4").Formula = "=A4*N4%"
$ws.Range("G4").Formula = "=A4*O4%"
$ws.Range("H4").Formula = "=SUM(B4:G4)"
$ws.Range("J4").Value = 24.1
$ws.Range("K4").Value = 3.7
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 22.2
$ws.Range("N4").Value = 3.7
$ws.Range("O4").Value = 46.3
$ws.Range("P4").Formula = "=SUM(J4:O4)"

# --- Row 5 ---
$ws.Range("A5").Value = 63
$ws.Range("B5").Formula = "=A5*J5%"
$ws.Range("C5").Formula = "=A5*K5%"
$ws.Range("D5").Formula = "=A5*L5%"
$ws.Range("E5").Formula = "=A5*M5%"
$ws.Range("F5").Formula = "=A5*N5%"
$ws.Range("G5").Formula = "=A5*O5%"
$ws.Range("H5").Formula = "=SUM(B5:G5)"
$ws.Range("J5").Value = 23.8
$ws.Range("K5").Value = 11.1
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 20.6
$ws.Range("N5").Value = 3.2
$ws.Range("O5").Value = 41.3
$ws.Range("P5").Formula = "=SUM(J5:O5)"

# --- Row 9: stray helper calc left behind by the author ---
$ws.Range("O9").Formula = "=(100-89.5)/2"

# --- Number formatting: integer format ("0" => numFmtId 1) on the computed
#     percentage-of-total columns B:I for rows 2-5 (I is blank but still
#     carries the format, matching the source columns being dragged across). ---
$ws.Range("B2:I5").NumberFormat = "0"

# --- Selection / active sheet: SiP becomes the active tab with N3 selected ---
$ws.Activate()
$ws.Range("N3").Select()
